# Add a new "PlacementStability" mapping row to the TestDataMappingSheet_SD
# sheet, right after the existing last row (row 151 / FolioInspections),
# mirroring the surrounding rows' layout:
#   A = ScreenName, B = TestDataFileName, C = TestDataSheetName,
#   D = StartIndexofIteration (1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 152

$ws.Cells.Item($newRow, 1).Value = "PlacementStability"
$ws.Cells.Item($newRow, 2).Value = "cares\Placement.xlsx"
$ws.Cells.Item($newRow, 3).Value = "PlacementStability"
$ws.Cells.Item($newRow, 4).Value = 1

# Match the selection/view left behind by the edit: both the new row and
# the row above it (D151:D152 were the last edited cells) end up selected.
$ws.Range("A151:D152").Select()

# Best-effort: scroll the window down so row 133 is the new top row
# (mirrors topLeftCell="A133" in the saved view).
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
